# raw_copy_expenses.xlsx update
#
# Context (see commit message):
#  - CSV import support was removed from the app; everything imported must
#    now come from an .xlsx file. The first few rows of this fixture held
#    data that used to be produced by the (now removed) CSV importer, so
#    those cells are cleared out.
#  - The date column shown by the notebook is now formatted without the
#    "T00:00:00" time component, i.e. "YYYY-MM-DDTHH:MM:SS" -> "YYYY-MM-DD".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the stale CSV-imported expense rows (rows 2-9) ---
# Each of these rows keeps its ID (column A) but loses the amount / currency
# / description / date / sum / category data that used to come from the CSV
# import.
$rowsToClear = @{
    2 = @("B2", "C2", "D2", "E2", "H2")
    3 = @("B3", "C3", "D3", "H3")
    4 = @("B4", "D4", "E4", "H4")
    5 = @("B5", "C5", "D5", "E5", "H5")
    6 = @("B6", "C6", "D6", "E6", "H6")
    7 = @("B7", "C7", "D7", "H7")
    8 = @("B8", "D8", "E8", "H8")
    9 = @("B9", "D9", "E9", "F9", "H9")
}

foreach ($rowNum in $rowsToClear.Keys) {
    foreach ($addr in $rowsToClear[$rowNum]) {
        $ws.Range($addr).Value = ""
    }
}

# --- 2. Reformat every remaining Date cell in column E ---
# Strip the trailing "T00:00:00" time portion, keeping the cell as plain
# text (not re-parsed into a date serial number).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$suffix = "T00:00:00"

for ($r = 10; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)   # column E = Date
    $v = $cell.Value2
    if ($v -eq $null) { continue }
    $s = [string]$v
    if ($s.EndsWith($suffix)) {
        $newValue = $s.Substring(0, $s.Length - $suffix.Length)
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
    }
}

# --- 3. Special case: row 10's date is corrected to 2019-08-17 ---
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2019-08-17"
